# Update the "想去人数" (want-to-go count) figures on the "展览" sheet
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value = 158
$wsExhibit.Range("F6").Value = 108
$wsExhibit.Range("F7").Value = 288
$wsExhibit.Range("F9").Value = 2031
$wsExhibit.Range("F11").Value = 4829
$wsExhibit.Range("F12").Value = 91

# Same underlying rows are duplicated on the "全部类型" (all-types) sheet,
# offset by two rows further down.
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 158
$wsAll.Range("F8").Value = 108
$wsAll.Range("F9").Value = 288
$wsAll.Range("F13").Value = 2031
$wsAll.Range("F15").Value = 4829
$wsAll.Range("F16").Value = 91
